$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '28.134.09'
Set-TextValue 'E2' '  +1.17%  '
Set-TextValue 'D3' '1.821.98'
Set-TextValue 'E3' '  +1.82%  '
Set-TextValue 'D4' '0.9992'
Set-TextValue 'E4' '  -0.19%  '
Set-TextValue 'D5' '311.18'
Set-TextValue 'E5' '  +0.35%  '
Set-TextValue 'D6' '0.9999'
Set-TextValue 'E6' '  -0.06%  '
Set-TextValue 'D7' '0.5016'
Set-TextValue 'E7' '  -2.12%  '
Set-TextValue 'D8' '0.3926'
Set-TextValue 'E8' '  +1.05%  '
Set-TextValue 'D9' '0.09921'
Set-TextValue 'E9' '  +26.97%  '
Set-TextValue 'D10' '1.110'
Set-TextValue 'E10' '  +1.81%  '
Set-TextValue 'D11' '40.90'
Set-TextValue 'E11' '  +0.28%  '
Set-TextValue 'D12' '6.432'
Set-TextValue 'E12' '  +3.53%  '
Set-TextValue 'D13' '20.62'
Set-TextValue 'E13' '  +2.19%  '
Set-TextValue 'D14' '0.9989'
Set-TextValue 'E14' '  -0.18%  '
Set-TextValue 'D15' '1.815.36'
Set-TextValue 'E15' '  +2.17%  '
Set-TextValue 'D16' '7.297'
Set-TextValue 'E16' '  +1.08%  '
Set-TextValue 'D17' '0.00001135'
Set-TextValue 'E17' '  +5.67%  '
Set-TextValue 'D18' '92.66'
Set-TextValue 'E18' '  +1.32%  '
Set-TextValue 'D19' '0.06652'
Set-TextValue 'E19' '  +2.04%  '
Set-TextValue 'D20' '0.9996'
Set-TextValue 'E20' '  -0.07%  '
Set-TextValue 'D21' '17.23'
Set-TextValue 'E21' '  +1.19%  '
Set-TextValue 'D22' '5.952'
Set-TextValue 'E22' '  +0.70%  '
Set-TextValue 'D23' '28.171.26'
Set-TextValue 'E23' '  +1.01%  '
Set-TextValue 'D24' '11.13'
Set-TextValue 'E24' '  +1.10%  '
Set-TextValue 'D25' '2.264'
Set-TextValue 'E25' '  +1.74%  '
Set-TextValue 'D26' '159.01'
Set-TextValue 'E26' '  -0.83%  '
Set-TextValue 'D27' '20.74'
Set-TextValue 'E27' '  +2.39%  '
Set-TextValue 'D28' '2.023.55'
Set-TextValue 'E28' '  +1.76%  '
Set-TextValue 'D29' '2.423'
Set-TextValue 'E29' '  +2.90%  '
Set-TextValue 'D30' '127.35'
Set-TextValue 'E30' '  +2.82%  '
Set-TextValue 'D31' '0.1068'
Set-TextValue 'E31' '  -0.62%  '
Set-TextValue 'D32' '1.041'
Set-TextValue 'E32' '  +0.41%  '
Set-TextValue 'D33' '5.582'
Set-TextValue 'E33' '  +1.78%  '
Set-TextValue 'D34' '3.599'
Set-TextValue 'E34' '  -0.59%  '
Set-TextValue 'D35' '0.06740'
Set-TextValue 'E35' '  -4.12%  '
Set-TextValue 'D39' '4.964'
Set-TextValue 'E39' '  -0.55%  '
Set-TextValue 'D40' '11.33'
Set-TextValue 'E40' '  -1.53%  '
Set-TextValue 'D41' '0.6209'
Set-TextValue 'E41' '  +2.07%  '
Set-TextValue 'D42' '1.177'
Set-TextValue 'E42' '  +2.52%  '
Set-TextValue 'D43' '0.9990'
Set-TextValue 'E43' '  -0.11%  '
Set-TextValue 'D44' '13.17'
Set-TextValue 'E44' '  +0.07%  '
Set-TextValue 'D45' '0.5926'
Set-TextValue 'E45' '  +0.61%  '
Set-TextValue 'D46' '3.697'
Set-TextValue 'E46' '  +0.07%  '
Set-TextValue 'D47' '1.282'
Set-TextValue 'E47' '  -2.19%  '
Set-TextValue 'D48' '124.53'
Set-TextValue 'E48' '  +0.94%  '
Set-TextValue 'D49' '1.937'
Set-TextValue 'E49' '  +1.49%  '
Set-TextValue 'D50' '1.185'
Set-TextValue 'E50' '  -1.20%  '
Set-TextValue 'D51' '0.06794'
Set-TextValue 'E51' '  -0.23%  '

# Row 36 and 37 swap (VeChain / FraxShare reordered with updated values)
Set-TextValue 'B36' 'VeChain'
Set-TextValue 'C36' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D36' '0.02345'
Set-TextValue 'E36' '  +1.84%  '
Set-TextValue 'B37' 'FraxShare'
Set-TextValue 'C37' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D37' '8.941'
Set-TextValue 'E37' '  +2.42%  '

# Row 38 - only Volume(1h) changes
Set-TextValue 'E38' '  +0.71%  '

